$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell C10 from 18 to 100 (value-only change per commit "SAVE")
$ws.Range("C10").Value = 100
